# TIMES_H2_Demand_PtL_2050.xlsx - "Correct excel files related to H2"
#
# Summary of the edit being reproduced:
#   - Header cell B1 text changes from "Process Description1" to "PJ" and is
#     restyled to match the B2/C2 header cells (Arial, centered, no
#     vertical-center).
#   - Column D (the "...LowAF" figures / "Supply convert CCS emission to
#     Methanol supplied by H2 LowAF" header) is removed entirely, shifting
#     everything left and shrinking the B1:D1 merge down to B1:C1.
#   - The sheet selection ends up on the merged header B1:C1.
#   - Page setup is switched to paper size 9 (A4) / portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Restyle B1 to match the header style already used by B2/C2
#    (Arial 666666, centered, quote-prefixed text) before changing its text.
$ws.Range("B2").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# 2. Change the B1 header text from "Process Description1" to "PJ".
#    Using a leading apostrophe keeps it as explicit text (quote-prefixed),
#    matching the original authoring.
$ws.Range("B1").Formula = "'PJ"

# 3. Remove column D completely (data, header "...LowAF" text, and column
#    width) - this shifts column C's neighbours left and drops the now
#    unused "Supply convert CCS emission to Methanol supplied by H2 LowAF"
#    string.
$ws.Columns("D:D").Delete()

# 4. The merged header cell is now B1:C1; select it, matching the saved
#    sheet view.
$ws.Range("B1:C1").Select()

# 5. Update the page setup (paper size 9 = A4, portrait orientation).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Applied H2 PtL 2050 corrections"
